$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.758.98"
$ws.Range("D3").Value = "3.141.67"
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'581.06"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("D6").Value = "'147.17"
$ws.Range("E6").Value = "  -2.51%  "
$ws.Range("D8").Value = "3.140.68"
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("D9").Value = "'0.526"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("D10").Value = "'0.157"
$ws.Range("E10").Value = "  -3.37%  "
$ws.Range("D11").Value = "'6.16"
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("E12").Value = "  -2.26%  "
$ws.Range("E13").Value = "  -3.03%  "
$ws.Range("D14").Value = "'37.09"
$ws.Range("E14").Value = "  -3.03%  "
$ws.Range("D15").Value = "3.656.20"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").Value = "64.835.68"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("E17").Value = "  -1.29%  "
$ws.Range("D18").Value = "3.139.58"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("D20").Value = "'499.31"
$ws.Range("D21").Value = "'15.40"
$ws.Range("E21").Value = "  +2.74%  "
$ws.Range("D22").Value = "'0.713"
$ws.Range("E22").Value = "  -3.57%  "
$ws.Range("D23").Value = "'14.99"
$ws.Range("E23").Value = "  -7.30%  "
$ws.Range("D24").Value = "'7.75"
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("D25").Value = "'84.19"
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").Value = "'9.10"
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("D30").Value = "'2.81"
$ws.Range("E30").Value = "  +1.17%  "
$ws.Range("D31").Value = "'27.56"
$ws.Range("E31").Value = "  -1.74%  "
$ws.Range("E32").Value = "  -0.76%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "'6.37"
$ws.Range("E34").Value = "  +1.07%  "
$ws.Range("E35").Value = "  -3.10%  "
$ws.Range("D36").Value = "'54.84"
$ws.Range("E36").Value = "  -1.84%  "
$ws.Range("D37").Value = "'0.0893"
$ws.Range("E37").Value = "  +1.63%  "
$ws.Range("D38").Value = "'469.33"
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("D39").Value = "'0.0417"
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("D40").Value = "'2.91"
$ws.Range("E40").Value = "  -6.69%  "
$ws.Range("D41").Value = "'8.74"
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("D42").Value = "2.978.28"
$ws.Range("E42").Value = "  -4.58%  "
$ws.Range("E43").Value = "  -4.03%  "
$ws.Range("D44").Value = "'2.42"
$ws.Range("E44").Value = "  -3.74%  "
$ws.Range("E45").Value = "  -3.22%  "
$ws.Range("D46").Value = "'28.19"
$ws.Range("E46").Value = "  -3.68%  "
$ws.Range("D47").Value = "0.0₃0597"
$ws.Range("E47").Value = "  +2.04%  "
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("E50").Value = "  -4.20%  "
$ws.Range("D51").Value = "'119.22"
$ws.Range("E51").Value = "  -4.27%  "
